$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column (B) holds plain "yyyy-mm-dd" strings. Force it to Text
# format before assigning so Excel doesn't auto-convert the literal
# string into a date serial number; the other columns are unambiguous
# text already, so they are left with their original formatting.
$ws.Range("B2:B3").NumberFormat = "@"

$ws.Range("A2").Value = "Chandy Neat"
$ws.Range("B2").Value = "2024-11-25"
$ws.Range("C2").Value = "07:39:25"
$ws.Range("D2").Value = "2024-11-25 07:39:25"

$ws.Range("A3").Value = "Koemthay Tha"
$ws.Range("B3").Value = "2024-11-25"
$ws.Range("C3").Value = "07:39:45"
$ws.Range("D3").Value = "2024-11-25 07:39:45"
